$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New row 21: same layout/format as row 6 (meeting row), except J21
#     takes the "OK-column" style (s=4 / K6's style) instead of F:I's (s=2) ---
$ws.Range("A6:K6").Copy()
$ws.Range("A21:K21").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("K6").Copy()
$ws.Range("J21").PasteSpecial(-4122)      # xlPasteFormats
$ws.Rows("21:21").RowHeight = 30

$ws.Range("A21").Value = "17/9/2014"
$ws.Range("B21").Value = "17/9/2014"
$ws.Range("C21").Value = "17/9/2014"
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = "17/9/2014"
$ws.Range("F21").Value = "Họp dự án"
$ws.Range("G21").Value = "*Tổng hợp kiến thức cá nhân mỗi người`n*Kế hoạch tiếp theo"
$ws.Range("H21").Value = "OK"
$ws.Range("I21").Value = "OK"
$ws.Range("J21").Value = "OK"
$ws.Range("K21").Value = "OK"

# --- New row 24: same layout/format as row 8, cell-by-cell so only the
#     occupied columns (A,B,C,D,F,G,I) get created, matching row 8's shape ---
$ws.Range("A8").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("B8").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$ws.Range("C8").Copy()
$ws.Range("C24").PasteSpecial(-4122)
$ws.Range("D8").Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("F8").Copy()
$ws.Range("F24").PasteSpecial(-4122)
$ws.Range("G8").Copy()
$ws.Range("G24").PasteSpecial(-4122)
$ws.Range("I8").Copy()
$ws.Range("I24").PasteSpecial(-4122)
$ws.Rows("24:24").RowHeight = 30

$ws.Range("A24").Value = "17/9/2014"
$ws.Range("F24").Value = "Xác nhận mail"
$ws.Range("G24").Value = "*Tìm hiểu xác nhận mail đăng nhập`n"
$ws.Range("I24").Value = "`n*Todo: Tìm hiểu về cơ chế sinh link, key để đăng ký acc cho website"
$ws.Range("B24").Value = "27/9/2014"

# --- Update sheet view: scroll position + selection ---
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("B24").Select()
